$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DC")

# Row 4 is being turned from the "q3" / "Vermogen" input question
# into a new (still mostly empty/placeholder) multiple-choice question.

$ws.Range("A4").Value = ""
$ws.Range("B4").Value = "mc"
$ws.Range("C4").Value = ""
$ws.Range("D4").Value = "Is dit een goede nieuwe vraag??"
$ws.Range("E4").Value = "['A. test 1', ' B. Test 2', ' C. Test 3']"
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = ""
$ws.Range("H4").Value = ""
$ws.Range("I4").Value = ""
$ws.Range("J4").Value = ""
$ws.Range("K4").Value = ""
